$d = $word.ActiveDocument

# 1) "City Data" table, "Size of the data" cell currently holds only "B".
#    A leading " k" was typed in front of it, so the cell reads " kB".
#    MatchWholeWord keeps this from touching the "B" inside "MB"/"kB" elsewhere.
$d.Content.Find.Execute("B", $true, $true, $false, $false, $false, $true, 1, $false, " kB", 2) | Out-Null

# 2) "Inflation Data" heading paragraph: the two runs ("Inflation" + " Data")
#    were retyped/merged into a single run with the same text.
$d.Content.Find.Execute("Inflation Data", $true, $false, $false, $false, $false, $true, 1, $false, "Inflation Data", 2) | Out-Null

# 3) "Inflation Data" table, "Size of the data" cell: runs "9.61k" + "B"
#    merged into a single run "9.61kB".
$d.Content.Find.Execute("9.61kB", $true, $false, $false, $false, $false, $true, 1, $false, "9.61kB", 2) | Out-Null

# 4) Second "Customer Data" table, "Total number of observations" cell:
#    runs "49" + "9" merged into a single run "499".
$d.Content.Find.Execute("499", $true, $false, $false, $false, $false, $true, 1, $false, "499", 2) | Out-Null
